$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "A101 (31 él.)"
$ws.Range("D2").Value = "D03 (87 él.)"
$ws.Range("E2").Value = "D03 (93 él.)"
$ws.Range("F2").Value = "A115 (14 él.)"
$ws.Range("B3").Value = "A116 (42 él.)"
$ws.Range("C3").Value = "A017 (43 él.)"
$ws.Range("D3").Value = "J108 (38 él.)"
$ws.Range("E3").Value = "J108 (38 él.)"
$ws.Range("F3").Value = "J020 (56 él.)"
$ws.Range("B4").Value = "J009 (13 él.)"
$ws.Range("C4").Value = "J108 (38 él.)"
$ws.Range("D4").Value = "J020 (125 él.)"
$ws.Range("E4").Value = "J012 (24 él.)"
$ws.Range("F4").Value = "J110 (23 él.)"
$ws.Range("B5").Value = "A017 (43 él.)"
$ws.Range("C5").Value = "J109 (38 él.)"
$ws.Range("D5").Value = "A115 (45 él.)"
$ws.Range("E5").Value = "J110 (43 él.)"
$ws.Range("F5").Value = "D03 (35 él.)"
$ws.Range("B6").Value = "J107 (13 él.)"
$ws.Range("C6").Value = "J107 (19 él.)"
$ws.Range("D6").Value = "J012 (24 él.)"
$ws.Range("E6").Value = "J022 (155 él.)"
$ws.Range("F6").Value = "-"
$ws.Range("B7").Value = "D03 (38 él.)"
$ws.Range("C7").Value = "D03 (96 él.)"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "A016 (37 él.)"
$ws.Range("F7").Value = "J109 (33 él.)"
$ws.Range("B8").Value = "J110 (36 él.)"
$ws.Range("C8").Value = "I013 (60 él.)"
$ws.Range("D8").Value = "A012 (38 él.)"
$ws.Range("E8").Value = "A116 (47 él.)"
$ws.Range("F8").Value = "J108 (30 él.)"
$ws.Range("B9").Value = "A012 (18 él.)"
$ws.Range("C9").Value = "J020 (134 él.)"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "J106 (36 él.)"
$ws.Range("F9").Value = "J204 (24 él.)"
$ws.Range("B10").Value = "A016 (19 él.)"
$ws.Range("C10").Value = "A016 (38 él.)"
$ws.Range("D10").Value = "I013 (60 él.)"
$ws.Range("E10").Value = "J020 (130 él.)"
$ws.Range("F10").Value = "-"
$ws.Range("B11").Value = "J108 (38 él.)"
$ws.Range("C11").Value = "A116 (48 él.)"
$ws.Range("D11").Value = "A013 (50 él.)"
$ws.Range("E11").Value = "I013 (60 él.)"
$ws.Range("F11").Value = "J022 (25 él.)"
$ws.Range("B12").Value = "J204 (20 él.)"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "J106 (36 él.)"
$ws.Range("E12").Value = "J021 (148 él.)"
$ws.Range("F12").Value = "J012 (17 él.)"
$ws.Range("B13").Value = "J020 (45 él.)"
$ws.Range("C13").Value = "A012 (38 él.)"
$ws.Range("D13").Value = "A116 (47 él.)"
$ws.Range("E13").Value = "J109 (38 él.)"
$ws.Range("F13").Value = "A116 (41 él.)"
$ws.Range("B14").Value = "A013 (41 él.)"
$ws.Range("C14").Value = "A013 (50 él.)"
$ws.Range("D14").Value = "A017 (43 él.)"
$ws.Range("E14").Value = "A115 (45 él.)"
$ws.Range("F14").Value = "A017 (37 él.)"
$ws.Range("B15").Value = "A115 (45 él.)"
$ws.Range("C15").Value = "J106 (36 él.)"
$ws.Range("D15").Value = "A016 (38 él.)"
$ws.Range("E15").Value = "A013 (50 él.)"
$ws.Range("F15").Value = "A012 (29 él.)"
$ws.Range("B16").Value = "J106 (12 él.)"
$ws.Range("C16").Value = "J022 (137 él.)"
$ws.Range("D16").Value = "A101 (31 él.)"
$ws.Range("E16").Value = "J204 (36 él.)"
$ws.Range("F16").Value = "J007 (9 él.)"
$ws.Range("B17").Value = "J022 (67 él.)"
$ws.Range("C17").Value = "J007 (19 él.)"
$ws.Range("D17").Value = "J110 (43 él.)"
$ws.Range("E17").Value = "A017 (43 él.)"
$ws.Range("F17").Value = "J021 (60 él.)"
$ws.Range("B18").Value = "I013 (33 él.)"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "J022 (153 él.)"
$ws.Range("E18").Value = "J008 (19 él.)"
$ws.Range("F18").Value = "A013 (20 él.)"
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "J021 (169 él.)"
$ws.Range("E19").Value = "A012 (38 él.)"
$ws.Range("F19").Value = "-"
$ws.Range("B20").Value = "J007 (7 él.)"
$ws.Range("C20").Value = "J021 (152 él.)"
$ws.Range("D20").Value = "J109 (32 él.)"
$ws.Range("E20").Value = "J107 (18 él.)"
$ws.Range("F20").Value = "-"
$ws.Range("B21").Value = "J109 (5 él.)"
$ws.Range("C21").Value = "J008 (19 él.)"
$ws.Range("D21").Value = "J107 (17 él.)"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
$ws.Range("B22").Value = "-"
$ws.Range("C22").Value = "-"
$ws.Range("D22").Value = "J204 (35 él.)"
$ws.Range("F22").Value = "J107 (5 él.)"
$ws.Range("B23").Value = "-"
$ws.Range("C23").Value = "J204 (34 él.)"
$ws.Range("D23").Value = "J008 (18 él.)"
$ws.Range("B24").Value = "J008 (15 él.)"
$ws.Range("D24").Value = "J007 (18 él.)"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "A101 (19 él.)"
$ws.Range("B25").Value = "J012 (18 él.)"
$ws.Range("C25").Value = "-"
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = "J106 (21 él.)"
$ws.Range("B26").Value = "A101 (6 él.)"
$ws.Range("D26").Value = "J009 (14 él.)"
$ws.Range("E26").Value = "J009 (18 él.)"
$ws.Range("C27").Value = "A115 (45 él.)"
$ws.Range("E27").Value = "-"
$ws.Range("F27").Value = "J009 (10 él.)"
$ws.Range("B28").Value = "-"
$ws.Range("C28").Value = "J012 (24 él.)"
$ws.Range("E28").Value = "A101 (26 él.)"
$ws.Range("B29").Value = "-"
$ws.Range("C29").Value = "J009 (19 él.)"
$ws.Range("F29").Value = "A016 (16 él.)"
$ws.Range("C30").Value = "-"
$ws.Range("E30").Value = "J007 (19 él.)"
$ws.Range("F30").Value = "I013 (22 él.)"
$ws.Range("B31").Value = "J021 (23 él.)"
$ws.Range("F31").Value = "J008 (18 él.)"
$ws.Range("C32").Value = "J110 (43 él.)"
$ws.Range("D32").Value = "-"
$ws.Range("F32").Value = "-"

Write-Host "Updated 134 cells"
